$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.381.08"
$ws.Range("E2").Value = "  +1.37%  "

$ws.Range("D3").Value = "1.858.06"
$ws.Range("E3").Value = "  +0.43%  "

$ws.Range("E4").Value = "  -0.23%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "310.99"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.44%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.010"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.19%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4769"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.29%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3801"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +3.41%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07301"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.98%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9292"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.33%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.73"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +5.06%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07793"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.58%  "

$ws.Range("D13").Value = "1.862.23"
$ws.Range("E13").Value = "  +0.96%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.439"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.75%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.537"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.51%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "89.98"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.05%  "

$ws.Range("E17").Value = "  -0.25%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008806"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.58%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.009"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.24%  "

$ws.Range("D20").Value = "27.448.02"
$ws.Range("E20").Value = "  +1.53%  "

$ws.Range("E21").Value = "  +0.81%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.093"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.37%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.68"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.33%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.944"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.58%  "

$ws.Range("E25").Value = "  +1.42%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.004"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.69%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "115.32"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.87%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "4.926"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.00%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.08882"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.24%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.325"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.59%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.208"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.24%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.588"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.84%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7488"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.89%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.712"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.56%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.122"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.67%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02036"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.93%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5529"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +4.84%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05255"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.22%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.987"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.27%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "7.025"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.19%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.590"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +4.00%  "

$ws.Range("E43").Value = "  +0.33%  "

$ws.Range("B44").Value = "Decentraland"
$ws.Range("C44").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.4872"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.48%  "

$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.60"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.34%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.011"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.24%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.658"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.98%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "102.77"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.03%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "67.35"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.63%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06075"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.05%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.9107"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.30%  "
